$wb = $excel.ActiveWorkbook

# --- Update Benchmark sheet values with the refreshed simulation results ---
$wsBench = $wb.Worksheets.Item("Benchmark")

$wsBench.Range("B2").Value = 174.18
$wsBench.Range("C2").Value = 148.63
$wsBench.Range("D2").Value = 117.58
$wsBench.Range("E2").Value = 180.82

$wsBench.Range("B4").Value = 24.74
$wsBench.Range("C4").Value = 18.46
$wsBench.Range("D4").Value = 14.48
$wsBench.Range("E4").Value = 22.65

$wsBench.Range("B5").Value = 174.18
$wsBench.Range("C5").Value = 148.63
$wsBench.Range("D5").Value = 117.58
$wsBench.Range("E5").Value = 180.82

$wsBench.Range("B6").Value = 2.36
$wsBench.Range("C6").Value = 1.65
$wsBench.Range("D6").Value = 1.12
$wsBench.Range("E6").Value = 3.63

$wsBench.Range("E7").Value = 9739.07

$wsBench.Range("B8").Value = 554
$wsBench.Range("C8").Value = 771
$wsBench.Range("D8").Value = 701
$wsBench.Range("E8").Value = 484

$wsBench.Range("B9").Value = -354
$wsBench.Range("C9").Value = -265
$wsBench.Range("D9").Value = -189
$wsBench.Range("E9").Value = -373

# --- Remove the "Time (Second)" sheet entirely ---
$excel.DisplayAlerts = $false
$wsTime = $wb.Worksheets.Item("Time (Second)")
$wsTime.Delete() | Out-Null
$excel.DisplayAlerts = $true

# --- Update sheet view selections on the remaining sheets ---
$wsCompetitive = $wb.Worksheets.Item("Competitive")
$wsCooperative = $wb.Worksheets.Item("Cooperative")

$wsCooperative.Range("D15").Select() | Out-Null
$wsBench.Range("B26").Select() | Out-Null

# "Competitive" becomes the active sheet/tab with cell C22 selected
$wsCompetitive.Activate() | Out-Null
$wsCompetitive.Range("C22").Select() | Out-Null
